# Regenerate save_data "K" column (column G) values.
# This mirrors the commit: "regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals" — here only the per-row K values
# (column G) change for this particular workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 0
    14 = 1
    15 = 2
    16 = 1
    17 = 0
    18 = 1
    19 = 2
    20 = 3
    22 = 1
    23 = 2
    24 = 3
    25 = 0
    26 = 0
    27 = 1
    28 = 1
    29 = 0
    30 = 2
    31 = 2
    33 = 0
    34 = 2
    35 = 0
    36 = 2
    37 = 1
    38 = 0
    39 = 2
    40 = 1
    41 = 1
    42 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
